$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-18 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-19 Monday", 2) | Out-Null
$d.Content.Find.Execute("16×99=", $true, $false, $false, $false, $false, $true, 1, $false, "91×13=", 2) | Out-Null
$d.Content.Find.Execute("79×75=", $true, $false, $false, $false, $false, $true, 1, $false, "34×77=", 2) | Out-Null
$d.Content.Find.Execute("16×69=", $true, $false, $false, $false, $false, $true, 1, $false, "64×36=", 2) | Out-Null
$d.Content.Find.Execute("14×13=", $true, $false, $false, $false, $false, $true, 1, $false, "23×75=", 2) | Out-Null
$d.Content.Find.Execute("97×89=", $true, $false, $false, $false, $false, $true, 1, $false, "21×21=", 2) | Out-Null
$d.Content.Find.Execute("22×58=", $true, $false, $false, $false, $false, $true, 1, $false, "24×50=", 2) | Out-Null
$d.Content.Find.Execute("18×42=", $true, $false, $false, $false, $false, $true, 1, $false, "44×57=", 2) | Out-Null
$d.Content.Find.Execute("32×16=", $true, $false, $false, $false, $false, $true, 1, $false, "78×37=", 2) | Out-Null
$d.Content.Find.Execute("32×83=", $true, $false, $false, $false, $false, $true, 1, $false, "81×78=", 2) | Out-Null
$d.Content.Find.Execute("24×88=", $true, $false, $false, $false, $false, $true, 1, $false, "79×33=", 2) | Out-Null
$d.Content.Find.Execute("29×31=", $true, $false, $false, $false, $false, $true, 1, $false, "37×67=", 2) | Out-Null
$d.Content.Find.Execute("11×59=", $true, $false, $false, $false, $false, $true, 1, $false, "66×95=", 2) | Out-Null
$d.Content.Find.Execute("51×86=", $true, $false, $false, $false, $false, $true, 1, $false, "19×34=", 2) | Out-Null
$d.Content.Find.Execute("78×97=", $true, $false, $false, $false, $false, $true, 1, $false, "38×19=", 2) | Out-Null
$d.Content.Find.Execute("65×26=", $true, $false, $false, $false, $false, $true, 1, $false, "76×18=", 2) | Out-Null
$d.Content.Find.Execute("97×80=", $true, $false, $false, $false, $false, $true, 1, $false, "38×20=", 2) | Out-Null
$d.Content.Find.Execute("29×73=", $true, $false, $false, $false, $false, $true, 1, $false, "35×39=", 2) | Out-Null
$d.Content.Find.Execute("68×17=", $true, $false, $false, $false, $false, $true, 1, $false, "47×29=", 2) | Out-Null
$d.Content.Find.Execute("89×43=", $true, $false, $false, $false, $false, $true, 1, $false, "53×75=", 2) | Out-Null
$d.Content.Find.Execute("81×81=", $true, $false, $false, $false, $false, $true, 1, $false, "18×52=", 2) | Out-Null
$d.Content.Find.Execute("23×40=", $true, $false, $false, $false, $false, $true, 1, $false, "39×54=", 2) | Out-Null
$d.Content.Find.Execute("44×33=", $true, $false, $false, $false, $false, $true, 1, $false, "89×27=", 2) | Out-Null
$d.Content.Find.Execute("74×20=", $true, $false, $false, $false, $false, $true, 1, $false, "38×100=", 2) | Out-Null
$d.Content.Find.Execute("80×99=", $true, $false, $false, $false, $false, $true, 1, $false, "72×98=", 2) | Out-Null
$d.Content.Find.Execute("100×28=", $true, $false, $false, $false, $false, $true, 1, $false, "26×15=", 2) | Out-Null
$d.Content.Find.Execute("45×97=", $true, $false, $false, $false, $false, $true, 1, $false, "10×73=", 2) | Out-Null
$d.Content.Find.Execute("38×76=", $true, $false, $false, $false, $false, $true, 1, $false, "57×73=", 2) | Out-Null
$d.Content.Find.Execute("64×55=", $true, $false, $false, $false, $false, $true, 1, $false, "81×25=", 2) | Out-Null
$d.Content.Find.Execute("66×85=", $true, $false, $false, $false, $false, $true, 1, $false, "66×67=", 2) | Out-Null
$d.Content.Find.Execute("80×24=", $true, $false, $false, $false, $false, $true, 1, $false, "70×16=", 2) | Out-Null
$d.Content.Find.Execute("25×16=", $true, $false, $false, $false, $false, $true, 1, $false, "30×33=", 2) | Out-Null
$d.Content.Find.Execute("26×57=", $true, $false, $false, $false, $false, $true, 1, $false, "96×98=", 2) | Out-Null
$d.Content.Find.Execute("61×79=", $true, $false, $false, $false, $false, $true, 1, $false, "83×52=", 2) | Out-Null
$d.Content.Find.Execute("81×56=", $true, $false, $false, $false, $false, $true, 1, $false, "40×37=", 2) | Out-Null
$d.Content.Find.Execute("62×19=", $true, $false, $false, $false, $false, $true, 1, $false, "36×86=", 2) | Out-Null
$d.Content.Find.Execute("28×100=", $true, $false, $false, $false, $false, $true, 1, $false, "28×56=", 2) | Out-Null
$d.Content.Find.Execute("54×66=", $true, $false, $false, $false, $false, $true, 1, $false, "47×52=", 2) | Out-Null
$d.Content.Find.Execute("39×15=", $true, $false, $false, $false, $false, $true, 1, $false, "78×65=", 2) | Out-Null
$d.Content.Find.Execute("97×50=", $true, $false, $false, $false, $false, $true, 1, $false, "82×62=", 2) | Out-Null
$d.Content.Find.Execute("79×24=", $true, $false, $false, $false, $false, $true, 1, $false, "48×19=", 2) | Out-Null
$d.Content.Find.Execute("87×21=", $true, $false, $false, $false, $false, $true, 1, $false, "13×44=", 2) | Out-Null
$d.Content.Find.Execute("41×62=", $true, $false, $false, $false, $false, $true, 1, $false, "56×42=", 2) | Out-Null
$d.Content.Find.Execute("30×45=", $true, $false, $false, $false, $false, $true, 1, $false, "57×53=", 2) | Out-Null
$d.Content.Find.Execute("93×57=", $true, $false, $false, $false, $false, $true, 1, $false, "93×95=", 2) | Out-Null
$d.Content.Find.Execute("67×50=", $true, $false, $false, $false, $false, $true, 1, $false, "90×16=", 2) | Out-Null
$d.Content.Find.Execute("86×48=", $true, $false, $false, $false, $false, $true, 1, $false, "88×74=", 2) | Out-Null
$d.Content.Find.Execute("25×65=", $true, $false, $false, $false, $false, $true, 1, $false, "11×23=", 2) | Out-Null
$d.Content.Find.Execute("77×66=", $true, $false, $false, $false, $false, $true, 1, $false, "87×25=", 2) | Out-Null
$d.Content.Find.Execute("57×26=", $true, $false, $false, $false, $false, $true, 1, $false, "64×97=", 2) | Out-Null
$d.Content.Find.Execute("29×78=", $true, $false, $false, $false, $false, $true, 1, $false, "24×47=", 2) | Out-Null
$d.Content.Find.Execute("18×14=", $true, $false, $false, $false, $false, $true, 1, $false, "57×61=", 2) | Out-Null
$d.Content.Find.Execute("19×98=", $true, $false, $false, $false, $false, $true, 1, $false, "34×100=", 2) | Out-Null
$d.Content.Find.Execute("27×10=", $true, $false, $false, $false, $false, $true, 1, $false, "86×90=", 2) | Out-Null
$d.Content.Find.Execute("12×15=", $true, $false, $false, $false, $false, $true, 1, $false, "44×39=", 2) | Out-Null
$d.Content.Find.Execute("15×64=", $true, $false, $false, $false, $false, $true, 1, $false, "40×36=", 2) | Out-Null
$d.Content.Find.Execute("67×93=", $true, $false, $false, $false, $false, $true, 1, $false, "18×77=", 2) | Out-Null
$d.Content.Find.Execute("46×41=", $true, $false, $false, $false, $false, $true, 1, $false, "62×18=", 2) | Out-Null
$d.Content.Find.Execute("54×71=", $true, $false, $false, $false, $false, $true, 1, $false, "19×32=", 2) | Out-Null
$d.Content.Find.Execute("47×49=", $true, $false, $false, $false, $false, $true, 1, $false, "69×26=", 2) | Out-Null
$d.Content.Find.Execute("89×80=", $true, $false, $false, $false, $false, $true, 1, $false, "67×64=", 2) | Out-Null
$d.Content.Find.Execute("48×71=", $true, $false, $false, $false, $false, $true, 1, $false, "38×47=", 2) | Out-Null
$d.Content.Find.Execute("69×41=", $true, $false, $false, $false, $false, $true, 1, $false, "19×35=", 2) | Out-Null
$d.Content.Find.Execute("99×89=", $true, $false, $false, $false, $false, $true, 1, $false, "70×65=", 2) | Out-Null
$d.Content.Find.Execute("13×26=", $true, $false, $false, $false, $false, $true, 1, $false, "80×98=", 2) | Out-Null
$d.Content.Find.Execute("90×50=", $true, $false, $false, $false, $false, $true, 1, $false, "85×59=", 2) | Out-Null
$d.Content.Find.Execute("93×51=", $true, $false, $false, $false, $false, $true, 1, $false, "39×39=", 2) | Out-Null
$d.Content.Find.Execute("98×33=", $true, $false, $false, $false, $false, $true, 1, $false, "66×87=", 2) | Out-Null
$d.Content.Find.Execute("17×36=", $true, $false, $false, $false, $false, $true, 1, $false, "32×71=", 2) | Out-Null
$d.Content.Find.Execute("94×77=", $true, $false, $false, $false, $false, $true, 1, $false, "70×69=", 2) | Out-Null
$d.Content.Find.Execute("66×90=", $true, $false, $false, $false, $false, $true, 1, $false, "87×72=", 2) | Out-Null
$d.Content.Find.Execute("19×65=", $true, $false, $false, $false, $false, $true, 1, $false, "41×42=", 2) | Out-Null
$d.Content.Find.Execute("29×64=", $true, $false, $false, $false, $false, $true, 1, $false, "87×19=", 2) | Out-Null
$d.Content.Find.Execute("72×21=", $true, $false, $false, $false, $false, $true, 1, $false, "72×76=", 2) | Out-Null
$d.Content.Find.Execute("59×37=", $true, $false, $false, $false, $false, $true, 1, $false, "70×14=", 2) | Out-Null
$d.Content.Find.Execute("51×65=", $true, $false, $false, $false, $false, $true, 1, $false, "41×96=", 2) | Out-Null
$d.Content.Find.Execute("34×14=", $true, $false, $false, $false, $false, $true, 1, $false, "54×10=", 2) | Out-Null
$d.Content.Find.Execute("31×68=", $true, $false, $false, $false, $false, $true, 1, $false, "89×90=", 2) | Out-Null
$d.Content.Find.Execute("76×62=", $true, $false, $false, $false, $false, $true, 1, $false, "66×54=", 2) | Out-Null
$d.Content.Find.Execute("83×89=", $true, $false, $false, $false, $false, $true, 1, $false, "44×63=", 2) | Out-Null
$d.Content.Find.Execute("24×60=", $true, $false, $false, $false, $false, $true, 1, $false, "19×21=", 2) | Out-Null
$d.Content.Find.Execute("65×77=", $true, $false, $false, $false, $false, $true, 1, $false, "10×75=", 2) | Out-Null
$d.Content.Find.Execute("75×84=", $true, $false, $false, $false, $false, $true, 1, $false, "83×51=", 2) | Out-Null
$d.Content.Find.Execute("62×99=", $true, $false, $false, $false, $false, $true, 1, $false, "68×29=", 2) | Out-Null
$d.Content.Find.Execute("48×28=", $true, $false, $false, $false, $false, $true, 1, $false, "82×47=", 2) | Out-Null
$d.Content.Find.Execute("88×72=", $true, $false, $false, $false, $false, $true, 1, $false, "35×47=", 2) | Out-Null
$d.Content.Find.Execute("18×13=", $true, $false, $false, $false, $false, $true, 1, $false, "49×23=", 2) | Out-Null
$d.Content.Find.Execute("69×30=", $true, $false, $false, $false, $false, $true, 1, $false, "88×25=", 2) | Out-Null
$d.Content.Find.Execute("22×23=", $true, $false, $false, $false, $false, $true, 1, $false, "55×85=", 2) | Out-Null
$d.Content.Find.Execute("47×13=", $true, $false, $false, $false, $false, $true, 1, $false, "71×64=", 2) | Out-Null
$d.Content.Find.Execute("75×13=", $true, $false, $false, $false, $false, $true, 1, $false, "38×73=", 2) | Out-Null
$d.Content.Find.Execute("50×73=", $true, $false, $false, $false, $false, $true, 1, $false, "50×75=", 2) | Out-Null
$d.Content.Find.Execute("56×95=", $true, $false, $false, $false, $false, $true, 1, $false, "61×58=", 2) | Out-Null
$d.Content.Find.Execute("23×82=", $true, $false, $false, $false, $false, $true, 1, $false, "30×26=", 2) | Out-Null
$d.Content.Find.Execute("46×37=", $true, $false, $false, $false, $false, $true, 1, $false, "46×83=", 2) | Out-Null
$d.Content.Find.Execute("92×61=", $true, $false, $false, $false, $false, $true, 1, $false, "44×35=", 2) | Out-Null
$d.Content.Find.Execute("98×66=", $true, $false, $false, $false, $false, $true, 1, $false, "52×16=", 2) | Out-Null
$d.Content.Find.Execute("83×69=", $true, $false, $false, $false, $false, $true, 1, $false, "58×81=", 2) | Out-Null
$d.Content.Find.Execute("23×90=", $true, $false, $false, $false, $false, $true, 1, $false, "38×69=", 2) | Out-Null
$d.Content.Find.Execute("48×93=", $true, $false, $false, $false, $false, $true, 1, $false, "10×12=", 2) | Out-Null
$d.Content.Find.Execute("40×40=", $true, $false, $false, $false, $false, $true, 1, $false, "63×26=", 2) | Out-Null
